# The underlying data rows for several observations got re-ordered/corrected.
# Rather than moving whole rows (which would risk Excel re-interpreting
# untouched date-like text cells such as the Startdatum/Slutdatum columns
# as real dates), we only swap the values of the specific cells that
# actually differ between each pair of rows. The row number itself never
# changes - only the data that lives in it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row A, row B, and the list of column letters whose values
# must be exchanged between row A and row B.
$swaps = @(
    @(2,  3,  @("A", "Q", "R", "Z", "AB")),
    @(6,  7,  @("A", "B", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB", "AC")),
    @(8,  9,  @("A", "B", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB", "AC", "AE")),
    @(10, 11, @("A", "B", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB", "AC")),
    @(14, 15, @("A", "B", "E", "F", "G", "H", "M", "Q", "R", "Z", "AB", "AC")),
    @(16, 17, @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")),
    @(25, 26, @("A", "Q", "R", "Z", "AB"))
)

foreach ($swap in $swaps) {
    $rowA = $swap[0]
    $rowB = $swap[1]
    $cols = $swap[2]

    foreach ($col in $cols) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
